$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.887.33'
$ws.Range("E2").Value = '  -1.59%  '
$ws.Range("D3").Value = '1.826.51'
$ws.Range("E3").Value = '  -1.61%  '
$ws.Range("E4").Value = '  +0.68%  '
$ws.Range("D5").Value = '''310.59'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("E6").Value = '  +0.57%  '
$ws.Range("D7").Value = '''0.4571'
$ws.Range("E7").Value = '  -0.86%  '
$ws.Range("D8").Value = '''0.3685'
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("D9").Value = '''0.07165'
$ws.Range("E9").Value = '  -2.35%  '
$ws.Range("D10").Value = '''0.8733'
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("D11").Value = '''0.07766'
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").Value = '''19.59'
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("D13").Value = '1.814.99'
$ws.Range("E13").Value = '  -2.60%  '
$ws.Range("D14").Value = '''5.315'
$ws.Range("E14").Value = '  -1.39%  '
$ws.Range("D15").Value = '''6.379'
$ws.Range("E15").Value = '  -2.56%  '
$ws.Range("E16").Value = '  -5.47%  '
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").Value = '''0.000008714'
$ws.Range("E18").Value = '  -3.34%  '
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").Value = '26.924.91'
$ws.Range("E20").Value = '  -1.53%  '
$ws.Range("E21").Value = '  -2.27%  '
$ws.Range("D22").Value = '''5.002'
$ws.Range("E22").Value = '  -2.34%  '
$ws.Range("D23").Value = '2.067.37'
$ws.Range("D24").Value = '''10.44'
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").Value = '''2.001'
$ws.Range("E25").Value = '  +4.51%  '
$ws.Range("D26").Value = '''151.54'
$ws.Range("D27").Value = '''18.17'
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").Value = '''1.956'
$ws.Range("E28").Value = '  -5.58%  '
$ws.Range("D29").Value = '''113.55'
$ws.Range("E29").Value = '  -2.22%  '
$ws.Range("D30").Value = '''4.900'
$ws.Range("E30").Value = '  -4.27%  '
$ws.Range("D31").Value = '''0.08794'
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("D32").Value = '''3.046'
$ws.Range("E32").Value = '  +1.45%  '
$ws.Range("E33").Value = '  -2.41%  '
$ws.Range("E34").Value = '  -0.27%  '
$ws.Range("E35").Value = '  -3.59%  '
$ws.Range("E36").Value = '  -4.12%  '
$ws.Range("D37").Value = '''1.086'
$ws.Range("E37").Value = '  +0.66%  '
$ws.Range("D38").Value = '''0.01945'
$ws.Range("E38").Value = '  -0.86%  '
$ws.Range("D39").Value = '''0.05132'
$ws.Range("E39").Value = '  -1.84%  '
$ws.Range("D40").Value = '''2.910'
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("D41").Value = '''6.934'
$ws.Range("E41").Value = '  -1.67%  '
$ws.Range("D42").Value = '''0.4968'
$ws.Range("E42").Value = '  -3.49%  '
$ws.Range("D43").Value = '''0.1595'
$ws.Range("E43").Value = '  -2.71%  '
$ws.Range("D44").Value = '''8.306'
$ws.Range("E44").Value = '  -0.71%  '
$ws.Range("D45").Value = '''0.4690'
$ws.Range("E45").Value = '  -3.05%  '
$ws.Range("D46").Value = '''1.006'
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("E47").Value = '  -1.38%  '
$ws.Range("D48").Value = '''101.57'
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("D49").Value = '''1.612'
$ws.Range("E49").Value = '  -2.53%  '
$ws.Range("D50").Value = '''0.06105'
$ws.Range("E50").Value = '  -1.81%  '
$ws.Range("D51").Value = '''64.50'
$ws.Range("E51").Value = '  -2.00%  '